# AFDP-7334 Admin Module - Holiday Schedule Misspelling
#
# Fixes the misspelled shared string "OCRWrokflow" -> "OCRWorkflow" used
# as the business-process name in the two rows of the "OCR Workflow
# Rules" decision table (column E, "Name of business process to start").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E17").Value = "OCRWorkflow"
$ws.Range("E18").Value = "OCRWorkflow"

# Leave the sheet with the same cell selected as in the authored edit.
$ws.Activate()
$ws.Range("D16").Select()
